# Add the AUDPC (Area Under Disease Progress Curve) Shapiro-Wilk normality
# test results as a new worksheet, mirroring the existing per-trait sheets.

$wb = $excel.ActiveWorkbook

# New sheet goes after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Lentil_FT13038_AUDPC_SW"

# Header row.
$newSheet.Range("A1").Value = "statistic"
$newSheet.Range("B1").Value = "p.value"
$newSheet.Range("C1").Value = "method"

# Shapiro-Wilk normality test result row for AUDPC.
$newSheet.Range("A2").Value = 0.986852229551412
$newSheet.Range("B2").Value = 0.200146571224564
$newSheet.Range("C2").Value = "Shapiro-Wilk normality test"

# Turn the range into a real Excel table, like the other sheets.
$tbl = $newSheet.ListObjects.Add(1, $newSheet.Range("A1:C2"), $null, 1)
$tbl.Name = "Table9"
$tbl.TableStyle = "TableStyleLight9"

# Restore the first sheet as the active tab/view (matches the committed
# workbook state: activeTab/firstSheet both reset to the first sheet).
$wb.Worksheets.Item(1).Activate()
